# Generate Report for Handoff
# Status moved from "In Translation" to "Ready for handoff" and the
# handoff timestamps were refreshed on the zh-cn and de-de report sheets.
# The "Status" columns are widened to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update the per-language status summary cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet: update status + latest handoff datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-29 15:08:50"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet: update status + latest handoff datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-29 15:08:55"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
